$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet and update the "through" date references
$ws.Name = "Through 2021-11-13"
$ws.Range("B1").Value = "November 2021 (through November 13)"

# Row 2 - North Lawndale
$ws.Range("B2").Value = 8
$ws.Range("M2").Value = 9

# Row 3 - Garfield Park
$ws.Range("AI3").Value = 2
$ws.Range("AT3").Value = 2
$ws.Range("BE3").Value = 3

# Row 5 - Humboldt Park
$ws.Range("AT5").Value = 1
$ws.Range("BP5").Value = 1

# Row 6 - West Town
$ws.Range("AT6").Value = 6

# Row 7 - South Shore
$ws.Range("M7").Value = 5

# Row 8 - Englewood
$ws.Range("AT8").Value = 1
$ws.Range("BE8").Value = 3

# Row 12 - Grand Boulevard
$ws.Range("B12").Value = 2

# Row 13 - Grand Crossing
$ws.Range("M13").Value = 1

# Row 22 - Avondale
$ws.Range("AI22").Value = 3
$ws.Range("AT22").Value = 1

# Row 23 - Lake View
$ws.Range("B23").Value = 2

# Row 44 - East Village
$ws.Range("AT44").Value = 2
$ws.Range("BE44").Value = 1

# Row 45 - United Center
$ws.Range("X45").Value = 1
$ws.Range("AT45").Value = 4

# Row 48 - Roseland
$ws.Range("AT48").Value = 3

# Row 54 - Hermosa
$ws.Range("M54").Value = 1
